$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.258.95"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "1.859.85"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7038"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "237.63"
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.08241"
$ws.Range("E8").Value = "  +10.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3040"
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.30"
$ws.Range("E10").Value = "  -0.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08188"
$ws.Range("E11").Value = "  +0.76%  "
$ws.Range("D12").Value = "1.895.10"
$ws.Range("E12").Value = "  +3.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7162"
$ws.Range("E13").Value = "  -1.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.179"
$ws.Range("E14").Value = "  -0.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.18"
$ws.Range("D16").Value = "29.279.97"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.779"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("E18").Value = "  +2.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "237.30"
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9999"
$ws.Range("D22").Value = "2.110.14"
$ws.Range("E22").Value = "  -0.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.460"
$ws.Range("E24").Value = "  -1.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.90"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.986"
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1443"
$ws.Range("E27").Value = "  -1.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.11"
$ws.Range("E28").Value = "  +0.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.966"
$ws.Range("E29").Value = "  +1.36%  "
$ws.Range("E30").Value = "  +4.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.425"
$ws.Range("E31").Value = "  -3.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.484"
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.061"
$ws.Range("E33").Value = "  +1.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05215"
$ws.Range("E34").Value = "  +1.16%  "
$ws.Range("E35").Value = "  -1.46%  "
$ws.Range("E36").Value = "  +1.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.001"
$ws.Range("E37").Value = "  -3.63%  "
$ws.Range("E38").Value = "  +0.91%  "
$ws.Range("E39").Value = "  -0.96%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.725"
$ws.Range("E40").Value = "  +1.74%  "
$ws.Range("D41").Value = "1.135.61"
$ws.Range("E41").Value = "  +5.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9162"
$ws.Range("E42").Value = "  -3.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.965"
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4288"
$ws.Range("E44").Value = "  -0.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.89"
$ws.Range("E45").Value = "  +1.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9994"
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.62"
$ws.Range("E47").Value = "  +0.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.774"
$ws.Range("E48").Value = "  +1.31%  "
$ws.Range("D49").Value = "2.008.05"
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.182"
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.975"
$ws.Range("E51").Value = "  -1.16%  "
